# Update TC_ID Excel SCD0017 until SCD0025 and Update TC_ID Solution SCD0006 until SCD0025

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab) from SCD0277 to SCD0017
$ws.Name = "SCD0017"

# Update TC_ID column (B) values for data rows 2-4 from "DGS-292" to "SCD0017-007"
$ws.Range("B2").Value = "SCD0017-007"
$ws.Range("B3").Value = "SCD0017-007"
$ws.Range("B4").Value = "SCD0017-007"

# Adjust column B width to fit the new, longer TC_ID text (matches Excel's bestFit recalculation)
$ws.Columns("B").ColumnWidth = 11.45

# Update the selection/view to reflect the user's final cursor position
$ws.Range("A1").Select()
$ws.Range("B5").Select()
